$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 182 (date 20-09-2021): add MOVE value and fix VIX value
$ws.Cells.Item(182, 2).Value = 59.58
$ws.Cells.Item(182, 3).Value = 25.71

# Append new rows for the following dates
$newRows = @(
    @{ Row = 183; Date = "21-09-2021"; MOVE = 59.2;  VIX = 24.36 },
    @{ Row = 184; Date = "22-09-2021"; MOVE = 55.92; VIX = 20.87 },
    @{ Row = 185; Date = "23-09-2021"; MOVE = 56.79; VIX = 18.63 },
    @{ Row = 186; Date = "24-09-2021"; MOVE = 58.46; VIX = 17.75 },
    @{ Row = 187; Date = "27-09-2021"; VIX = 18.37 }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Date
    if ($item.ContainsKey("MOVE")) {
        $ws.Cells.Item($r, 2).Value = $item.MOVE
    }
    $ws.Cells.Item($r, 3).Value = $item.VIX
}
